# Update the "想去人数" (want-to-go count) figures for the two data rows
# on both the "展览" sheet and the "全部类型" sheet, matching the refreshed
# gh-pages data export (commit 456a3b4):
#   F2: 436 -> 438
#   F3: 11  -> 12

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 438
    $ws.Range("F3").Value = 12
}
